# Renames the inline logo pictures embedded in the document's headers and
# footers. The Pearson Edexcel logo (alt text / description ends in
# "PearsonLogo.png") is renamed from image1.png to image2.png, and the BTEC
# logo (alt text "BTec_Logo-Orange") is renamed from image2.jpg to
# image1.jpg, wherever they occur across the document's sections.

$d = $word.ActiveDocument

function Rename-LogoInlineShapes($range) {
    $count = $range.InlineShapes.Count
    for ($k = 1; $k -le $count; $k++) {
        $shape = $range.InlineShapes.Item($k)
        $desc = $shape.AlternativeText

        if ($desc -like "*PearsonLogo.png") {
            $shape.Name = "image2.png"
        } elseif ($desc -eq "BTec_Logo-Orange") {
            $shape.Name = "image1.jpg"
        }
    }
}

for ($si = 1; $si -le $d.Sections.Count; $si++) {
    $section = $d.Sections.Item($si)

    $headers = $section.Headers
    for ($hi = 1; $hi -le $headers.Count; $hi++) {
        $header = $headers.Item($hi)
        if ($header.Exists) {
            Rename-LogoInlineShapes $header.Range
        }
    }

    $footers = $section.Footers
    for ($fi = 1; $fi -le $footers.Count; $fi++) {
        $footer = $footers.Item($fi)
        if ($footer.Exists) {
            Rename-LogoInlineShapes $footer.Range
        }
    }
}
